$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.475.44"
$ws.Range("E2").Value = "  -1.09%  "

# Row 3
$ws.Range("D3").Value = "1.656.41"
$ws.Range("E3").Value = "  -2.87%  "

# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "

# Row 5
$ws.Range("D5").Value = "'307.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "

# Row 7
$ws.Range("D7").Value = "'0.3612"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.44%  "

# Row 8
$ws.Range("D8").Value = "'47.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.59%  "

# Row 9
$ws.Range("D9").Value = "'0.3241"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.86%  "

# Row 10
$ws.Range("D10").Value = "'1.117"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.12%  "

# Row 11
$ws.Range("D11").Value = "'0.06988"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.36%  "

# Row 12
$ws.Range("E12").Value = "  +0.53%  "

# Row 13
$ws.Range("D13").Value = "'5.870"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.00%  "

# Row 14
$ws.Range("D14").Value = "'19.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.50%  "

# Row 15
$ws.Range("D15").Value = "1.654.46"
$ws.Range("E15").Value = "  -2.98%  "

# Row 16
$ws.Range("D16").Value = "'6.551"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.83%  "

# Row 17
$ws.Range("D17").Value = "'0.00001043"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.31%  "

# Row 18
$ws.Range("D18").Value = "'0.06525"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.89%  "

# Row 19
$ws.Range("E19").Value = "  +0.28%  "

# Row 20
$ws.Range("D20").Value = "'76.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.26%  "

# Row 21
$ws.Range("D21").Value = "'5.905"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.77%  "

# Row 22
$ws.Range("D22").Value = "'15.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.75%  "

# Row 23
$ws.Range("D23").Value = "'12.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.33%  "

# Row 24
$ws.Range("D24").Value = "24.467.21"
$ws.Range("E24").Value = "  -0.96%  "

# Row 25
$ws.Range("E25").Value = "  +1.71%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.312"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -16.66%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'146.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.35%  "

# Row 28
$ws.Range("D28").Value = "'18.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.94%  "

# Row 29
$ws.Range("D29").Value = "1.840.87"
$ws.Range("E29").Value = "  -2.76%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'123.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.63%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.182"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "

# Row 32
$ws.Range("D32").Value = "'3.953"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.30%  "

# Row 33
$ws.Range("D33").Value = "'5.627"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -16.98%  "

# Row 34
$ws.Range("D34").Value = "'1.697"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.95%  "

# Row 35
$ws.Range("D35").Value = "'0.08370"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.02%  "

# Row 36
$ws.Range("D36").Value = "'12.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.84%  "

# Row 37
$ws.Range("D37").Value = "'5.180"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.31%  "

# Row 38
$ws.Range("D38").Value = "'0.06043"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.89%  "

# Row 39
$ws.Range("D39").Value = "'0.02198"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.78%  "

# Row 40
$ws.Range("D40").Value = "'1.200"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.94%  "

# Row 41
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.2047"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.36%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'8.185"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.68%  "

# Row 43
$ws.Range("E43").Value = "  +0.46%  "

# Row 44
$ws.Range("D44").Value = "'0.5892"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.26%  "

# Row 45
$ws.Range("D45").Value = "'3.735"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "

# Row 46
$ws.Range("D46").Value = "'12.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.54%  "

# Row 47
$ws.Range("D47").Value = "'0.5574"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.71%  "

# Row 48
$ws.Range("D48").Value = "'121.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.26%  "

# Row 49
$ws.Range("D49").Value = "'1.928"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.79%  "

# Row 50
$ws.Range("D50").Value = "'0.06897"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.28%  "

# Row 51
$ws.Range("D51").Value = "'74.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.51%  "
